$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# SamplesTab (row 3) query: drop the Tumor / Analyte Type columns from the
# SELECT list (the rest of the query is unchanged).
$ws.Range("B3").Value = @"
SELECT
    DISTINCT (smp.sample_id) AS "Sample ID",
    sp.participant_id AS "Participant ID", 
    s.study_name AS "Study Name",
    s.phs_accession AS Accession
FROM 
    df_participant sp
JOIN 
    df_study s ON sp."study.phs_accession" = s.phs_accession
JOIN 
    df_sample smp ON smp."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_diagnosis d ON d."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_program p ON p.program_acronym = s."program.program_acronym"
JOIN
    df_file f1 ON f1."sample.sample_id" = smp.sample_id
JOIN
    df_genomic_info gi ON gi."file.file_id" = f1.file_id
WHERE 
    s.phs_accession = 'phs001437' AND f1.experimental_strategy_and_data_subtypes = 'RNA-Seq|WXS'
ORDER BY 
    smp.sample_id ASC
LIMIT 100;
"@

# The TsvExcel / WebExcel helper columns (D/E) are no longer needed for the
# SamplesTab and FilesTab rows.
$ws.Range("D3").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("D4").ClearContents()
$ws.Range("E4").ClearContents()

# Move the active selection onto the SamplesTab query cell.
$ws.Range("C3").Select()
